$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" cells hold numeric-looking text (e.g. "233.00", "42.071.95")
# that must round-trip as literal text rather than being reinterpreted by
# Excel as a number (which would drop trailing zeros / reformat). Force the
# cells that need it to Text format before writing their values.
$textCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D46", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '42.071.95'
$ws.Range('E2').Value = '  -4.23%  '
$ws.Range('D3').Value = '2.239.83'
$ws.Range('E3').Value = '  -4.86%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '233.00'
$ws.Range('E5').Value = '  -3.30%  '
$ws.Range('D6').Value = '0.635'
$ws.Range('E6').Value = '  -6.17%  '
$ws.Range('D7').Value = '70.08'
$ws.Range('E7').Value = '  -4.69%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').Value = '0.561'
$ws.Range('E9').Value = '  -7.42%  '
$ws.Range('D10').Value = '0.0992'
$ws.Range('E10').Value = '  -1.58%  '
$ws.Range('D11').Value = '58.30'
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('D12').Value = '35.44'
$ws.Range('E12').Value = '  +5.21%  '
$ws.Range('E13').Value = '  -3.01%  '
$ws.Range('D14').Value = '6.82'
$ws.Range('E14').Value = '  -7.14%  '
$ws.Range('D15').Value = '2.573.83'
$ws.Range('E15').Value = '  -4.83%  '
$ws.Range('D16').Value = '15.03'
$ws.Range('E16').Value = '  -8.38%  '
$ws.Range('D17').Value = '0.862'
$ws.Range('E17').Value = '  -5.21%  '
$ws.Range('D18').Value = '2.242.07'
$ws.Range('E18').Value = '  -4.76%  '
$ws.Range('D19').Value = '42.007.13'
$ws.Range('E19').Value = '  -4.22%  '
$ws.Range('D20').Value = '0.0₃0982'
$ws.Range('E20').Value = '  -4.66%  '
$ws.Range('D21').Value = '6.25'
$ws.Range('E21').Value = '  -7.21%  '
$ws.Range('D22').Value = '73.33'
$ws.Range('E22').Value = '  -5.82%  '
$ws.Range('D23').Value = '237.32'
$ws.Range('E23').Value = '  -7.26%  '
$ws.Range('D24').Value = '2.05'
$ws.Range('E24').Value = '  +5.50%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  -2.72%  '
$ws.Range('D27').Value = '2.36'
$ws.Range('E27').Value = '  -6.14%  '
$ws.Range('D28').Value = '10.05'
$ws.Range('E28').Value = '  -5.77%  '
$ws.Range('D29').Value = '2.15'
$ws.Range('E29').Value = '  -5.60%  '
$ws.Range('D30').Value = '167.92'
$ws.Range('D31').Value = '20.65'
$ws.Range('E31').Value = '  -9.15%  '
$ws.Range('D32').Value = '0.119'
$ws.Range('E32').Value = '  -7.58%  '
$ws.Range('E33').Value = '  -7.21%  '
$ws.Range('D34').Value = '5.42'
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('D35').Value = '0.0717'
$ws.Range('E35').Value = '  -5.39%  '
$ws.Range('D36').Value = '4.78'
$ws.Range('E36').Value = '  -8.35%  '
$ws.Range('D37').Value = '3.60'
$ws.Range('E37').Value = '  -5.79%  '
$ws.Range('D38').Value = '22.16'
$ws.Range('E38').Value = '  +16.02%  '
$ws.Range('D39').Value = '6.06'
$ws.Range('E39').Value = '  -6.65%  '
$ws.Range('E40').Value = '  -5.80%  '
$ws.Range('D41').Value = '67.72'
$ws.Range('E41').Value = '  -0.83%  '
$ws.Range('D42').Value = '0.0267'
$ws.Range('E42').Value = '  -4.16%  '
$ws.Range('D43').Value = '9.11'
$ws.Range('E43').Value = '  -1.80%  '
$ws.Range('D44').Value = '4.92'
$ws.Range('E44').Value = '  -4.84%  '
$ws.Range('E45').Value = '  -10.25%  '
$ws.Range('D46').Value = '0.189'
$ws.Range('E46').Value = '  -7.18%  '
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '2.36'
$ws.Range('E48').Value = '  -5.68%  '
$ws.Range('B49').Value = 'SynthetixNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D49').Value = '4.37'
$ws.Range('E49').Value = '  +5.85%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').Value = '1.18'
$ws.Range('E50').Value = '  -6.87%  '
$ws.Range('B51').Value = 'Celestia'
$ws.Range('C51').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D51').Value = '10.08'
$ws.Range('E51').Value = '  +5.18%  '
